# Actualización automática hashcode jue abr  9 01:41:04 CEST 2020
# Update hashcode values (column B) for several rows in the metadata sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "B44"  = "a32cd81b35935c925f43d51a1956ce3f"
    "B89"  = "540c6e9b1efc86a7027d6bfbd80c73c2"
    "B99"  = "3ed806b97270274a88c3d0a88769021f"
    "B110" = "1cbee20c6dd597308e23e402c1cb3429"
    "B154" = "7883f0f152cc9d9bb5a1fc710f211227"
    "B159" = "aaac34bf96dd1a6394dd3ca31665d0c9"
    "B160" = "86c3466b53645a70143a60d23010a457"
    "B168" = "a332483bc4229b143d1abf32232184ec"
    "B222" = "60e39cb58668f837f3ef9ef35b8dd94c"
    "B229" = "7d4adc924049e1e26eb3f440c3450a2b"
    "B278" = "ff0cdaad1bb498b10fd0b974320bdfa6"
    "B281" = "d47b4c2c37695aeaedf46052fc07213c"
    "B335" = "ce0d246ac8e46bde9469712017fd6d68"
    "B339" = "0cfcf0cdbc873d2da6b6d2d79315cafe"
    "B523" = "c85280c7cb5f69f7fdc4117e7b066ac0"
    "B542" = "b526e2e952a95b9a09ec2a8738f95769"
    "B561" = "5cbb749084cfb11e073fabbd9fa5cca4"
    "B592" = "2a0370be441331729a17ae4b1bdd77b2"
    "B776" = "ec7cbf44da2741d451e3a0d8eb8e7bff"
    "B819" = "19e459ae140fd3ca9c68c0372a062362"
    "B823" = "ce02acf55c77ea096712c1a555e3035c"
    "B827" = "af8a0fdf3300e2447c7ee9846c20357a"
    "B833" = "138c1287037ebf103f817fe612d3f27d"
    "B835" = "820a409f29375b7c62388a0b687f0f64"
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
